$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "date" header to "QC event date"
$ws.Range("D1").Value = "QC event date"

# Remove the file*_type columns (H = file1_type, J = file2_type, L = file3_type).
# Delete from rightmost to leftmost so the remaining column letters stay valid.
$ws.Columns("L").Delete()
$ws.Columns("J").Delete()
$ws.Columns("H").Delete()

# Add two new attached-file columns (file4, file5) after the existing file3 column
$ws.Range("J1").Value = "file4"
$ws.Range("K1").Value = "file5"

# Move the active selection to D1, matching the saved view state
$ws.Range("D1").Select()
